$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","H","J","K","L","M","O")

$data = @(
    @(0.7117080123049107, 0.1600642082812556, 0.1049859070985804, 0.1422608133583871, 2.307636691385639, 0.07973214163530429, 0.213893841753837, 0.3364172470826645, 0.2358674811643695, 0.2029199323890225, 6.28772141877289),
    @(0.680430168657324, 0.1594244198751653, 0.1037915478909213, 0.1427241575656346, 2.317594672145361, 0.07973214163530429, 0.2151499831879065, 0.307317597068149, 0.2341666889501184, 0.1974650344182436, 6.321571816392634),
    @(0.6614864130268074, 0.159033878832993, 0.1030930743046383, 0.1430468559831102, 2.324613667606371, 0.07973214163530429, 0.2159684523110137, 0.2895201197245143, 0.2332019661925244, 0.1941925677695799, 6.344846736306749),
    @(0.6538329321449794, 0.1588753248073331, 0.1028172593191954, 0.1431879886222767, 2.327701745878962, 0.07973214163530429, 0.2163138712992367, 0.282285527889627, 0.2328289283031992, 0.192878473115865, 6.354958314533178),
    @(0.652566097621019, 0.158849033489691, 0.1027719945762158, 0.1432120059459514, 2.328228284332774, 0.07973214163530429, 0.2163719463158422, 0.2810853318373177, 0.2327682018288826, 0.1926614481192424, 6.356675211878326),
    @(0.6613829264105107, 0.1590317380860569, 0.1030893187930388, 0.1430487203221205, 2.324654391874304, 0.07973214163530429, 0.2159730626003604, 0.2894224778882659, 0.2331968537971605, 0.1941747664497484, 6.344980565708994),
    @(0.7008696152831533, 0.1598431443660786, 0.1045668851639476, 0.1424126596263875, 2.31088258294001, 0.07973214163530429, 0.2143171758988842, 0.3263694656214682, 0.2352645836105012, 0.2010232090928312, 6.298876441222035),
    @(0.7803502673314142, 0.1614517753212112, 0.1077390394662174, 0.1414674535086036, 2.29104387585371, 0.07973214163530429, 0.2114435662720506, 0.3993596803011883, 0.2399473750426822, 0.2150578027648891, 6.228205355032486),
    @(0.8399681022684717, 0.1626434518755886, 0.1102346882391956, 0.1409558612950885, 2.280824541167505, 0.07973214163530429, 0.2095588293747035, 0.4532961433035894, 0.243766694545954, 0.2257321894802686, 6.188288291082756),
    @(0.8673506231747581, 0.1631875420343789, 0.1114053863086184, 0.1407625532655388, 2.277118650276378, 0.07973214163530429, 0.2087503334748515, 0.4778974195051262, 0.2455856696066974, 0.2306660087951826, 6.172729985771383),
    @(0.8777568094151036, 0.1633938447698782, 0.1118537450049644, 0.1406949975765013, 2.275850674554775, 0.07973214163530429, 0.2084511867428773, 0.4872222785722045, 0.2462861162184424, 0.2325454139542416, 6.167211865722436),
    @(0.8755140151274929, 0.1633494021790156, 0.1117569595739241, 0.1407092961569241, 2.276117739063366, 0.07973214163530429, 0.2085153016840575, 0.485213615930121, 0.2461347461754571, 0.2321401600885622, 6.168383688126056),
    @(0.8682060082526846, 0.1632045094284038, 0.1114421723214818, 0.1407568823946335, 2.27701162154672, 0.07973214163530429, 0.2087255820298992, 0.4786644059785488, 0.2456430629762352, 0.2308204074305067, 6.172268524218879),
    @(0.8637344443854147, 0.1631157927350628, 0.1112500109565886, 0.1407867649351573, 2.277576772090285, 0.07973214163530429, 0.2088552976611755, 0.4746539686142057, 0.2453434061671587, 0.2300134589564138, 6.174696722946464),
    @(0.83818380671363, 0.1626079327620289, 0.1101588885683356, 0.1409692857030684, 2.281085684610886, 0.07973214163530429, 0.2096126487256083, 0.4516896633125782, 0.2436494538185414, 0.2254113101461215, 6.18935734829455),
    @(0.822575966303134, 0.1622968736169241, 0.1094985539795985, 0.1410913361947692, 2.283479634752723, 0.07973214163530429, 0.2100897672085829, 0.4376181755848734, 0.2426310937847376, 0.2226079158844456, 6.199016827186227),
    @(0.8136234561353888, 0.1621181490831773, 0.1091220835035571, 0.1411652479356107, 2.284945338755364, 0.07973214163530429, 0.2103687949291242, 0.4295308050772917, 0.2420530408016504, 0.2210028265078066, 6.204817469531605),
    @(0.8105965537377244, 0.1620576689646782, 0.1089951917011405, 0.1411909113221554, 2.285456853760948, 0.07973214163530429, 0.2104640597670091, 0.4267936366061349, 0.2418586438032833, 0.2204606372563447, 6.206823519164914),
    @(0.8242348953202736, 0.1623299670249878, 0.1095685027835103, 0.141077959744921, 2.283215609253915, 0.07973214163530429, 0.2100385009788148, 0.4391154738114551, 0.242738705643319, 0.2229055824229533, 6.197963229949409),
    @(0.8703515476067594, 0.1632470608255971, 0.1115344965765601, 0.140742752113896, 2.276745394651599, 0.07973214163530429, 0.2086636273867093, 0.4805878304240423, 0.2457871670237566, 0.2312077513200137, 6.171117321160864),
    @(0.9007069181457439, 0.1638479859152611, 0.112848744507815, 0.1405565726510023, 2.273305638846779, 0.07973214163530429, 0.207805938593518, 0.5077439969854822, 0.2478473210973675, 0.2366981868353051, 6.155748653429129),
    @(0.8844861892951883, 0.1635271250502441, 0.1121446365598615, 0.1406529373274505, 2.275069391904708, 0.07973214163530429, 0.2082599687742821, 0.4932456852502582, 0.2467416020528006, 0.2337619834413545, 6.163752175038553),
    @(0.8234848291547507, 0.1623150051594422, 0.10953686906754, 0.1410839955763521, 2.283334696718683, 0.07973214163530429, 0.2100616637188306, 0.4384385374116846, 0.2426900312169025, 0.2227709867522094, 6.198438790960978),
    @(0.7586320231601462, 0.1610148149816695, 0.1068517296984197, 0.1416909511064759, 2.295644810000844, 0.07973214163530429, 0.2121810882802269, 0.3795581356592379, 0.2386137366506631, 0.2111968977021874, 6.245213508256541)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowNum = $r + 2
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $addr = "$($cols[$c])$rowNum"
        $ws.Range($addr).Value = $data[$r][$c]
    }
}
